$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALZ Policy Assignments H1CY24")
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19,9).Value = "Subnets should be private - 7bca8353-aa3b-429b-904a-9229c4385837 (azadvertizer.net)"
$hl = $ws.Hyperlinks.Add($ws.Cells.Item(19,9), "https://www.azadvertizer.net/azpolicyadvertizer/7bca8353-aa3b-429b-904a-9229c4385837.html")
$ws.Cells.Item(19,9).WrapText = $true

$ws.Rows.Item(47).Insert()
$ws.Cells.Item(47,9).Value = "Subnets should be private - 7bca8353-aa3b-429b-904a-9229c4385837 (azadvertizer.net)"
$hl2 = $ws.Hyperlinks.Add($ws.Cells.Item(47,9), "https://www.azadvertizer.net/azpolicyadvertizer/7bca8353-aa3b-429b-904a-9229c4385837.html")
$ws.Cells.Item(47,9).WrapText = $true

Write-Output "Style I19:"
Write-Output $ws.Cells.Item(19,9).Style.Name
Write-Output "Style I47:"
Write-Output $ws.Cells.Item(47,9).Style.Name
